$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded replacement file drops the "Grayson Allen" row and re-enters
# the remaining players/positions/teams in a new order. Remove the now
# unused last row (19) so the table shrinks back to 18 rows (A1:C18), then
# rewrite the data rows (2-18) with the final values.
$ws.Rows(19).Delete()

$ws.Range("A2").Value = "Donovan Mitchell"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Cleveland Cavaliers"

$ws.Range("A3").Value = "Cam Thomas"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Brooklyn Nets"

$ws.Range("A4").Value = "Aaron Wiggins"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "Oklahoma City Thunder"

$ws.Range("A5").Value = "Malik Beasley"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Detroit Pistons"

$ws.Range("A6").Value = "De'Andre Hunter"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Cleveland Cavaliers"

$ws.Range("A7").Value = "Michael Porter Jr."
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Denver Nuggets"

$ws.Range("A8").Value = "Jaden McDaniels"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Minnesota Timberwolves"

$ws.Range("A9").Value = "Dyson Daniels"
$ws.Range("B9").Value = "PG,SG,SF"
$ws.Range("C9").Value = "Atlanta Hawks"

$ws.Range("A10").Value = "Kristaps Porzingis"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Boston Celtics"

$ws.Range("A11").Value = "Toumani Camara"
$ws.Range("B11").Value = "SG,SF,PF"
$ws.Range("C11").Value = "Portland Trail Blazers"

$ws.Range("A12").Value = "Josh Hart"
$ws.Range("B12").Value = "SG,SF,PF"
$ws.Range("C12").Value = "New York Knicks"

$ws.Range("A13").Value = "Alperen Sengün"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Houston Rockets"

$ws.Range("A14").Value = "Domantas Sabonis"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = "Sacramento Kings"

$ws.Range("A15").Value = "Kelly Oubre Jr."
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Philadelphia 76ers"

$ws.Range("A16").Value = "Devin Vassell"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "San Antonio Spurs"

$ws.Range("A17").Value = "Julius Randle"
$ws.Range("B17").Value = "PF,C"
$ws.Range("C17").Value = "Minnesota Timberwolves"

$ws.Range("A18").Value = "Guerschon Yabusele"
$ws.Range("B18").Value = "PF,C"
$ws.Range("C18").Value = "Philadelphia 76ers"
